$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 0.1511493333333333
$ws.Range("H2").Value = 0.453448
$ws.Range("I2").Value = 0.7495144539818079
$ws.Range("J2").Value = 0.7495144539818078
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 3.847811666666666
$ws.Range("N2").Value = 11.543435
$ws.Range("O2").Value = 0.0396810199351781
$ws.Range("P2").Value = 0.03968101993517809
$ws.Range("Q2").Value = 0.5815941682088889
$ws.Range("R2").Value = 5.234347513879999
$ws.Range("S2").Value = 0.02974149799015625
$ws.Range("T2").Value = 0.02974149799015624
# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 0.1511493333333333
$ws.Range("H3").Value = 0.453448
$ws.Range("I3").Value = 0.7495144539818079
$ws.Range("J3").Value = 0.7495144539818078
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 57.80210733333333
$ws.Range("N3").Value = 173.406322
$ws.Range("O3").Value = 0.5960911739155557
$ws.Range("P3").Value = 0.5960911739155557
$ws.Range("Q3").Value = 8.736749988695111
$ws.Range("R3").Value = 78.63074989825598
$ws.Range("S3").Value = 0.4467789507406926
$ws.Range("T3").Value = 0.4467789507406926
# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 0.1511493333333333
$ws.Range("H4").Value = 0.453448
$ws.Range("I4").Value = 0.7495144539818079
$ws.Range("J4").Value = 0.7495144539818078
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 6.066157
$ws.Range("N4").Value = 18.198471
$ws.Range("O4").Value = 0.06255797260873913
$ws.Range("P4").Value = 0.06255797260873913
$ws.Range("Q4").Value = 0.9168955864453334
$ws.Range("R4").Value = 8.252060278008001
$ws.Range("S4").Value = 0.04688810468204801
$ws.Range("T4").Value = 0.046888104682048
# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 0.1511493333333333
$ws.Range("H5").Value = 0.453448
$ws.Range("I5").Value = 0.7495144539818079
$ws.Range("J5").Value = 0.7495144539818078
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 29.25249166666667
$ws.Range("N5").Value = 87.757475
$ws.Range("O5").Value = 0.3016698335405271
$ws.Range("P5").Value = 0.301669833540527
$ws.Range("Q5").Value = 4.421494613755556
$ws.Range("R5").Value = 39.79345152379999
$ws.Range("S5").Value = 0.226105900568911
$ws.Range("T5").Value = 0.226105900568911
# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 0.05051366666666667
$ws.Range("H6").Value = 0.151541
$ws.Range("I6").Value = 0.2504855460181921
$ws.Range("J6").Value = 0.2504855460181921
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 3.847811666666666
$ws.Range("N6").Value = 11.543435
$ws.Range("O6").Value = 0.0396810199351781
$ws.Range("P6").Value = 0.03968101993517809
$ws.Range("Q6").Value = 0.1943670759261111
$ws.Range("R6").Value = 1.749303683335
$ws.Range("S6").Value = 0.00993952194502185
$ws.Range("T6").Value = 0.009939521945021849
# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 0.05051366666666667
$ws.Range("H7").Value = 0.151541
$ws.Range("I7").Value = 0.2504855460181921
$ws.Range("J7").Value = 0.2504855460181921
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 57.80210733333333
$ws.Range("N7").Value = 173.406322
$ws.Range("O7").Value = 0.5960911739155557
$ws.Range("P7").Value = 0.5960911739155557
$ws.Range("Q7").Value = 2.919796382466889
$ws.Range("R7").Value = 26.278167442202
$ws.Range("S7").Value = 0.149312223174863
$ws.Range("T7").Value = 0.149312223174863
# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 0.05051366666666667
$ws.Range("H8").Value = 0.151541
$ws.Range("I8").Value = 0.2504855460181921
$ws.Range("J8").Value = 0.2504855460181921
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 6.066157
$ws.Range("N8").Value = 18.198471
$ws.Range("O8").Value = 0.06255797260873913
$ws.Range("P8").Value = 0.06255797260873913
$ws.Range("Q8").Value = 0.3064238326456667
$ws.Range("R8").Value = 2.757814493811
$ws.Range("S8").Value = 0.01566986792669112
$ws.Range("T8").Value = 0.01566986792669112
# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 0.05051366666666667
$ws.Range("H9").Value = 0.151541
$ws.Range("I9").Value = 0.2504855460181921
$ws.Range("J9").Value = 0.2504855460181921
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 29.25249166666667
$ws.Range("N9").Value = 87.757475
$ws.Range("O9").Value = 0.3016698335405271
$ws.Range("P9").Value = 0.301669833540527
$ws.Range("Q9").Value = 1.477650613219445
$ws.Range("R9").Value = 13.298855518975
$ws.Range("S9").Value = 0.07556393297161604
$ws.Range("T9").Value = 0.07556393297161602
